$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (F2:AO2)
$ws.Range("F2").Value = 1.58
$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 4.8
$ws.Range("I2").Value = 7.8
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.2
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.83
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.32
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.84
$ws.Range("U2").Value = 1.84
$ws.Range("V2").Value = 1.14
$ws.Range("W2").Value = 2.1
$ws.Range("X2").Value = 18
$ws.Range("Y2").Value = 22
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 8.800000000000001
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 28
$ws.Range("AE2").Value = 110
$ws.Range("AF2").Value = 11
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 24
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 48
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 13
$ws.Range("AO2").Value = 1000

# Row 3 updates
$ws.Range("F3").Value = 1.97
$ws.Range("Q3").Value = 1.61
